$wb = $excel.ActiveWorkbook

# Physical sheet order is fixed: Worksheets.Item(1) == xl/worksheets/sheet1.xml,
# Worksheets.Item(2) == xl/worksheets/sheet2.xml. We swap what logically lives on
# each physical sheet: sheet1.xml becomes "review_info" (headers only), and
# sheet2.xml becomes "hotel_info" (headers + the one data row, with a new
# "State" column inserted after "Hotel_Name").
$wsA = $wb.Worksheets.Item(1)
$wsB = $wb.Worksheets.Item(2)

# Avoid name collisions while renaming.
$wsA.Name = "tmp_shard83_a"
$wsB.Name = "tmp_shard83_b"

# Wipe existing content/formatting so the rebuilt shared-strings table reflects
# only what we write below (and in the order we write it).
$wsA.Cells.Clear()
$wsB.Cells.Clear()

# ---- sheet1.xml -> "review_info" (header row only, no data rows) ----
$wsA.Name = "review_info"

$reviewHeaders = @("STR","reviewer_ID","reviewer_name","Review_ID","Date_of_scraping","ReviewURL","Tripadvisor_gcode","Tripadvisor_dcode","Tripadvisor_rcode","review_date","review_title","review_content","review_rating","trip_month","trip_purpose","value","rooms","Location","Cleanliness","Sleep Quality","Service","Picture(yes=1)","respondent","response_date","response_text")
for ($i = 0; $i -lt $reviewHeaders.Count; $i++) {
    $wsA.Cells.Item(1, $i + 1).Value = $reviewHeaders[$i]
}

# ---- sheet2.xml -> "hotel_info" (header row + one data row, with new State column) ----
$wsB.Name = "hotel_info"

$hotelHeaders = @("STR","Hotel_Name","State","City","Zip","TA_ReviewURL","Tripadvisor_Hotel_Name","English_Reviews_num","Local_Rank","Total_Reviews_num")
for ($i = 0; $i -lt $hotelHeaders.Count; $i++) {
    $wsB.Cells.Item(1, $i + 1).Value = $hotelHeaders[$i]
}

$wsB.Cells.Item(2, 1).Value = 40813
$wsB.Cells.Item(2, 2).Value = "Hilton Garden Inn New Orleans Convention Center"
$wsB.Cells.Item(2, 3).Value = "Louisiana"
$wsB.Cells.Item(2, 4).Value = "New Orleans"
$wsB.Cells.Item(2, 5).Value = 70130
$wsB.Cells.Item(2, 6).Value = "https://www.tripadvisor.com/Hotel_Review-g60864-d223120-Reviews-Hilton_Garden_Inn_New_Orleans_Convention_Center-New_Orleans_Louisiana.html"
$wsB.Cells.Item(2, 7).Value = "Hilton Garden Inn New Orleans Convention Center"

# English_Reviews_num / Local_Rank / Total_Reviews_num are numeric-looking but
# stored as TEXT in the source data, so force text format before assigning,
# then drop back to the Normal style so no stray style index lingers on the cell.
foreach ($pair in @(@(8, "1362"), @(9, "117"), @(10, "1408"))) {
    $col = $pair[0]
    $text = $pair[1]
    $cell = $wsB.Cells.Item(2, $col)
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.Style = "Normal"
}
